$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3909
$ws.Range("I3").Value = 4052
$ws.Range("I4").Value = 946
$ws.Range("I5").Value = 374
$ws.Range("I6").Value = 4513
$ws.Range("I7").Value = 13794

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I2").Value = 24
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I3").Value = 142
$ws.Range("I7").Value = 443

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 68
$ws.Range("I7").Value = 257

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 127
$ws.Range("I3").Value = 191
$ws.Range("I5").Value = 15
$ws.Range("I6").Value = 177
$ws.Range("I7").Value = 536

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I2").Value = 42
$ws.Range("I3").Value = 35
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 100
$ws.Range("I3").Value = 85
$ws.Range("I5").Value = 15
$ws.Range("I7").Value = 308

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 119
$ws.Range("I7").Value = 432
$ws.Range("I8").Value = 831
$ws.Range("I9").Value = 61
$ws.Range("I11").Value = 211
$ws.Range("I14").Value = 75
$ws.Range("I15").Value = 163
$ws.Range("I18").Value = 94
$ws.Range("I19").Value = 379
$ws.Range("I21").Value = 76
$ws.Range("I23").Value = 135
$ws.Range("I24").Value = 35
$ws.Range("I25").Value = 67
$ws.Range("I27").Value = 129
$ws.Range("I29").Value = 890
$ws.Range("I31").Value = 135
$ws.Range("I33").Value = 628
$ws.Range("I34").Value = 65
$ws.Range("I36").Value = 194
$ws.Range("I37").Value = 443
$ws.Range("I41").Value = 61
$ws.Range("I42").Value = 479
$ws.Range("I43").Value = 119
$ws.Range("I51").Value = 134
$ws.Range("I52").Value = 302
$ws.Range("I53").Value = 149
$ws.Range("I63").Value = 51
$ws.Range("I64").Value = 122
$ws.Range("I65").Value = 308
$ws.Range("I66").Value = 40
$ws.Range("I67").Value = 536
$ws.Range("I76").Value = 209
$ws.Range("I77").Value = 75
$ws.Range("I78").Value = 201
$ws.Range("I79").Value = 372
$ws.Range("I81").Value = 13
$ws.Range("I83").Value = 278
$ws.Range("I85").Value = 623
$ws.Range("I88").Value = 123
$ws.Range("I91").Value = 172
$ws.Range("I95").Value = 220
$ws.Range("I98").Value = 91
$ws.Range("I99").Value = 257
$ws.Range("I101").Value = 13794

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 278

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I3").Value = 84
$ws.Range("I7").Value = 220

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 233
$ws.Range("I5").Value = 25
$ws.Range("I6").Value = 195
$ws.Range("I7").Value = 628

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 305
$ws.Range("I4").Value = 44
$ws.Range("I6").Value = 242
$ws.Range("I7").Value = 890

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 108
$ws.Range("I6").Value = 104
$ws.Range("I7").Value = 379

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I2").Value = 44
$ws.Range("I4").Value = 27
$ws.Range("I7").Value = 209

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 161
$ws.Range("I3").Value = 251
$ws.Range("I6").Value = 154
$ws.Range("I7").Value = 623

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I3").Value = 21
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I6").Value = 130
$ws.Range("I7").Value = 479

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I3").Value = 52
$ws.Range("I7").Value = 201

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("I3").Value = 17
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I2").Value = 54
$ws.Range("I3").Value = 60
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 172

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I2").Value = 6
$ws.Range("I6").Value = 58
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 118
$ws.Range("I6").Value = 111
$ws.Range("I7").Value = 372

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 122

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I2").Value = 28
$ws.Range("I6").Value = 39
$ws.Range("I7").Value = 94

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 61
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 194

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 85
$ws.Range("I3").Value = 104
$ws.Range("I6").Value = 75
$ws.Range("I7").Value = 302

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I2").Value = 24
$ws.Range("I7").Value = 65

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I2").Value = 52
$ws.Range("I3").Value = 37
$ws.Range("I7").Value = 163

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 58
$ws.Range("I7").Value = 91

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("I2").Value = 12
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 95
$ws.Range("I7").Value = 211

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("I3").Value = 22
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 38
$ws.Range("I7").Value = 119

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I3").Value = 44
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 259
$ws.Range("I3").Value = 232
$ws.Range("I6").Value = 269
$ws.Range("I7").Value = 831

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("I4").Value = 18
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I4").Value = 14
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 134

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I2").Value = 23
$ws.Range("I7").Value = 119

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I6").Value = 68
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 75

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 131
$ws.Range("I7").Value = 432

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("I2").Value = 7
$ws.Range("I6").Value = 13
